# Auto-generated script applying scheduled market-data refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2258.077
$ws.Range("I17").Value = 2500
$ws.Range("J17").Value = 2106.875
$ws.Range("K17").Value = 7500
$ws.Range("L17").Value = 6320.625
$ws.Range("M17").Value = -7332
$ws.Range("N17").Value = -6656.625
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53736
$ws.Range("H116").Value = 2242.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2242.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 2242.5
$ws.Range("N116").Value = -9126.5
$ws.Range("H137").Value = 1114225.1
$ws.Range("I137").Value = 1114225.1
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3342675.3
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3340125.3
$ws.Range("H138").Value = 250
$ws.Range("I138").Value = 250
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 750
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 4390

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").ClearContents()
$ws.Range("H2").Value = 3100
$ws.Range("I2").Value = 3100
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3100
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = -2987
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 150
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -38
$ws.Range("H32").Value = 867
$ws.Range("I32").Value = 867
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 867
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -580
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1788
$ws.Range("H63").Value = 5689.6
$ws.Range("I63").Value = 5316.6665
$ws.Range("J63").Value = 6249
$ws.Range("K63").Value = 5316.6665
$ws.Range("L63").Value = 6249
$ws.Range("M63").Value = -4630.6665
$ws.Range("N63").Value = -7621
$ws.Range("H66").Value = 5689.6
$ws.Range("I66").Value = 5316.6665
$ws.Range("J66").Value = 6249
$ws.Range("K66").Value = 26583.3325
$ws.Range("L66").Value = 31245
$ws.Range("M66").Value = -23151.3325
$ws.Range("N66").Value = -38109
$ws.Range("M116").ClearContents()
$ws.Range("H116").Value = 3100
$ws.Range("I116").Value = 3100
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3100
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = -806
$ws.Range("M132").ClearContents()
$ws.Range("H132").Value = 9090
$ws.Range("I132").Value = 9090
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 27270
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -24740
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").ClearContents()
$ws.Range("H3").Value = 3100
$ws.Range("I3").Value = 3100
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3100
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = -2986
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -35
$ws.Range("H86").Value = 709
$ws.Range("I86").Value = 714.2857
$ws.Range("J86").Value = 699.75
$ws.Range("K86").Value = 714.2857
$ws.Range("L86").Value = 699.75
$ws.Range("M86").Value = 408.7143
$ws.Range("N86").Value = -2945.75
$ws.Range("H89").Value = 709
$ws.Range("I89").Value = 714.2857
$ws.Range("J89").Value = 699.75
$ws.Range("K89").Value = 3571.4285
$ws.Range("L89").Value = 3498.75
$ws.Range("M89").Value = 2044.5715
$ws.Range("N89").Value = -14730.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1149.6666
$ws.Range("I31").Value = 1149.6666
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1149.6666
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -854.6666
$ws.Range("H34").Value = 1149.6666
$ws.Range("I34").Value = 1149.6666
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1149.6666
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -947.6666
$ws.Range("H88").Value = 16510.545
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 16510.545
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 16510.545
$ws.Range("N88").Value = -17322.545
$ws.Range("H91").Value = 16510.545
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 16510.545
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 16510.545
$ws.Range("N91").Value = -19318.545

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4997.5
$ws.Range("I4").Value = 8000
$ws.Range("J4").Value = 1995
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 5985
$ws.Range("M4").Value = -23888
$ws.Range("N4").Value = -6209
$ws.Range("L23").ClearContents()
$ws.Range("H23").Value = 499.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 499.5
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = 1498.5
$ws.Range("N23").Value = -1968.5
$ws.Range("M86").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2247.6667
$ws.Range("I80").Value = 2121.75
$ws.Range("J80").Value = 2499.5
$ws.Range("K80").Value = 2121.75
$ws.Range("L80").Value = 2499.5
$ws.Range("M80").Value = -1123.75
$ws.Range("N80").Value = -4495.5
$ws.Range("H83").Value = 2247.6667
$ws.Range("I83").Value = 2121.75
$ws.Range("J83").Value = 2499.5
$ws.Range("K83").Value = 10608.75
$ws.Range("L83").Value = 12497.5
$ws.Range("M83").Value = -5616.75
$ws.Range("N83").Value = -22481.5
$ws.Range("H107").Value = 1655
$ws.Range("I107").Value = 1950
$ws.Range("J107").Value = 475
$ws.Range("K107").Value = 1950
$ws.Range("L107").Value = 475
$ws.Range("M107").Value = -30
$ws.Range("N107").Value = -4315
$ws.Range("H132").Value = 669333.3
$ws.Range("I132").Value = 669333.3
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2007999.9
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2005469.9
$ws.Range("H136").Value = 24999
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 24999
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 74997
$ws.Range("N136").Value = -80097

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4500
$ws.Range("I61").Value = 4500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4298
$ws.Range("H93").Value = 1439.2727
$ws.Range("I93").Value = 854
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 854
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = 394
$ws.Range("N93").Value = -5496
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2330
